$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A79").Value = "94EXRO"
$ws.Range("B79").Value = "Led con base 6 volt"
$ws.Range("D79").Value = 2500
$ws.Range("E79").Value = 7000
$ws.Range("F79").Value = 84
$ws.Range("G79").Value = 16
$ws.Range("H79").Formula = "=(E79-D79)*G79"
$ws.Range("I79").Formula = "=D79*F79"
$ws.Range("J79").Value = 210000
